# This workbook holds the naive QoQ forecast-error matrix (rows = forecast
# origin dates, columns B:K = Q0..Q9 horizons) matched against the ifoCAST
# series. The commit "added ifoCAST full series evaluation" rolls every row
# forward one quarter: column B (the oldest horizon) is dropped, everything
# shifts left by one column, and rows that were still accruing history (2-10)
# gain one freshly observed data point in column K. Rows that had already run
# out of columns (11-20) simply lose their trailing (rightmost) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K")

# Row number -> new left-to-right values for columns B.. (one entry per column).
$rowData = [ordered]@{
    2 = @(-1.579792844692939, -0.9155274120541215, -0.1460550576300293, -0.2386646425575917, 0.4822622889884105, -0.2266896356254442, 0.2654163566970329, 0.73522110397002, 1.048046298935466, 0.169643717683886)
    3 = @(-3.063097883853514, -2.293625529429422, -2.386235114356984, -1.665308182810982, -2.374260107424837, -1.88215411510236, -1.412349367829372, -1.099524172863926, -1.977926754115506, -2.707838769293975)
    4 = @(1.146404457093908, 1.053794872166346, 1.774721803712348, 1.065769879098493, 1.55787587142097, 2.027680618693958, 2.340505813659403, 1.462103232407824, 0.7321912172293545, 2.394370313618982)
    5 = @(1.114969050580547, 1.835895982126549, 1.126944057512694, 1.619050049835171, 2.088854797108159, 2.401679992073604, 1.523277410822025, 0.7933653956435556, 2.455544492033183, 1.836325203637939)
    6 = @(-0.2055599550297054, -0.9145118796435601, -0.422405887321083, 0.04739885995190407, 0.36022405491735, -0.5181785263342299, -1.248090541512699, 0.4140885548769285, -0.2051307335183153, 0.03217303010139827)
    7 = @(-0.3319173426099191, 0.160188649712558, 0.6299933969855451, 0.942818591950991, 0.06441601069941108, -0.6654960044790579, 0.9966830919105695, 0.3774638035153257, 0.6147675671350392, 0.5114901972596275)
    8 = @(0.0645136089032002, 0.5343183561761873, 0.8471435511416332, -0.03125903010994671, -0.7611710452884157, 0.9010080511012117, 0.2817887627059679, 0.5190925263256815, 0.4158151564502698, -0.3757832708791649)
    9 = @(0.7947373931749101, 1.107562588140356, 0.2291600068887761, -0.5007520082896928, 1.161427088099934, 0.5422077997046907, 0.7795115633244043, 0.6762341934489926, -0.1153642338804421, 0.4951246737870189)
    10 = @(0.2721075919818648, -0.6062949892697151, -1.336207004448184, 0.3259720919414433, -0.2932471964538005, -0.05594343283408693, -0.1592208027094986, -0.9508192300389333, -0.3403303223714723, -0.5018523531907899)
    11 = @(-0.5062156293670936, -1.236127644545562, 0.4260514518440648, -0.193167836551179, 0.04413592706853459, -0.05914144280687711, -0.8507398701363118, -0.2402509624688508, -0.4017729932881683)
    12 = @(-0.8823411384658664, 0.779837957923761, 0.1606186695285172, 0.3979224331482308, 0.2946450632728191, -0.4969533640566156, 0.1135355436108454, -0.04798648720847212)
    13 = @(0.6316605674913157, 0.0124412790960719, 0.2497450427157855, 0.1464676728403738, -0.6451307544890609, -0.03464184682159993, -0.1961638776409175)
    14 = @(-0.1450047099080831, 0.0922990537116305, -0.0109783161637812, -0.8025767434932158, -0.1920878358257549, -0.3536098666450724)
    15 = @(0.5929585102377013, 0.4896811403622896, -0.3019172869671451, 0.3085716207003159, 0.1470495898809984)
    16 = @(0.2093232598268204, -0.5822751675026142, 0.02821374016484672, -0.1333082906544708)
    17 = @(-0.4103003096576026, 0.2001885980098584, 0.03866656719054083)
    18 = @(-0.1056739417364731, -0.2671959725557906)
    19 = @(0.3451339801314955)
    20 = @()
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $addr = "$($cols[$i])$r"
        if ($i -lt $vals.Count) {
            $ws.Range($addr).Value = $vals[$i]
        } else {
            $ws.Range($addr).ClearContents()
        }
    }
}
